$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an explicit text value without leaving a lingering
# NumberFormat / style change on the cell (mirrors the original plain inline-string cells).
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "66.323.73"
Set-TextValue $ws.Range("E2") "  -1.42%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.563.00"
Set-TextValue $ws.Range("E3") "  +1.34%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.03%  "

# Row 5
Set-TextValue $ws.Range("D5") "609.16"
Set-TextValue $ws.Range("E5") "  -0.34%  "

# Row 6
Set-TextValue $ws.Range("D6") "144.56"
Set-TextValue $ws.Range("E6") "  -2.54%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.564.45"
Set-TextValue $ws.Range("E7") "  +1.57%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -0.03%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +0.29%  "

# Row 10
Set-TextValue $ws.Range("E10") "  -3.88%  "

# Row 11
Set-TextValue $ws.Range("D11") "8.09"
Set-TextValue $ws.Range("E11") "  +0.47%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.411"
Set-TextValue $ws.Range("E12") "  -2.69%  "

# Row 13
Set-TextValue $ws.Range("D13") "4.166.40"
Set-TextValue $ws.Range("E13") "  +1.30%  "

# Row 14
Set-TextValue $ws.Range("E14") "  -3.99%  "

# Row 15
Set-TextValue $ws.Range("D15") "30.19"
Set-TextValue $ws.Range("E15") "  -4.34%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.558.12"
Set-TextValue $ws.Range("E16") "  +1.12%  "

# Row 17
Set-TextValue $ws.Range("D17") "66.397.62"
Set-TextValue $ws.Range("E17") "  -1.40%  "

# Row 18
Set-TextValue $ws.Range("E18") "  -1.04%  "

# Row 19
Set-TextValue $ws.Range("D19") "11.36"
Set-TextValue $ws.Range("E19") "  +4.14%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -2.24%  "

# Row 21
Set-TextValue $ws.Range("D21") "14.98"
Set-TextValue $ws.Range("E21") "  -2.92%  "

# Row 22
Set-TextValue $ws.Range("D22") "429.41"
Set-TextValue $ws.Range("E22") "  -1.65%  "

# Row 23
Set-TextValue $ws.Range("E23") "  -0.71%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -1.52%  "

# Row 25
Set-TextValue $ws.Range("D25") "3.700.25"
Set-TextValue $ws.Range("E25") "  +1.16%  "

# Row 26
Set-TextValue $ws.Range("E26") "  -0.03%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0000122"
Set-TextValue $ws.Range("E27") "  +2.73%  "

# Row 28
Set-TextValue $ws.Range("D28") "8.12"
Set-TextValue $ws.Range("E28") "  -1.93%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -6.45%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -1.33%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.00"
Set-TextValue $ws.Range("E31") "  +0.08%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -5.52%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.159"
Set-TextValue $ws.Range("E33") "  -3.91%  "

# Row 34
Set-TextValue $ws.Range("D34") "25.45"
Set-TextValue $ws.Range("E34") "  -0.51%  "

# Row 35
Set-TextValue $ws.Range("D35") "3.552.33"
Set-TextValue $ws.Range("E35") "  +1.22%  "

# Row 36
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D36") "1.00"
Set-TextValue $ws.Range("E36") "  -0.03%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "1.76"
Set-TextValue $ws.Range("E37") "  -3.03%  "

# Row 38
Set-TextValue $ws.Range("D38") "7.84"
Set-TextValue $ws.Range("E38") "  -2.38%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.65"
Set-TextValue $ws.Range("E39") "  -5.39%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.02%  "

# Row 41
Set-TextValue $ws.Range("D41") "174.74"
Set-TextValue $ws.Range("E41") "  -0.91%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -2.34%  "

# Row 44
Set-TextValue $ws.Range("E44") "  +0.03%  "

# Row 45
Set-TextValue $ws.Range("D45") "1.91"
Set-TextValue $ws.Range("E45") "  -6.97%  "

# Row 46
Set-TextValue $ws.Range("D46") "45.70"
Set-TextValue $ws.Range("E46") "  -1.44%  "

# Row 47
Set-TextValue $ws.Range("E47") "  -1.41%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -7.60%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -2.30%  "

# Row 50
Set-TextValue $ws.Range("E50") "  -4.37%  "

# Row 51
Set-TextValue $ws.Range("D51") "23.05"
Set-TextValue $ws.Range("E51") "  +6.52%  "
